$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each data value is written with a leading apostrophe so Excel treats
# numeric-looking strings (e.g. '1.011') as literal text instead of
# silently parsing them into a Double. ClearFormats() then strips the
# quote-prefix marker Excel stamps on the cell style for that, restoring
# the cell's original (unstyled) look.

$ws.Range("D2").Value = "'" + '27.454.24'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'" + '  +1.72%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'" + '1.862.28'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'" + '  +0.72%  '
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'" + '1.011'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'" + '  -0.26%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'" + '311.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'" + '  +0.38%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'" + '1.010'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'" + '0.4775'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'" + '  -0.15%  '
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'" + '0.3802'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'" + '  +3.31%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'" + '0.07329'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'" + '  +1.32%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'" + '0.9350'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'" + '  +0.69%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'" + '20.75'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'" + '  +5.24%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'" + '0.07808'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'" + '  +0.85%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'" + '1.886.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'" + '  +3.18%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'" + '5.441'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'" + '  +1.75%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'" + '6.552'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'" + '  +1.77%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'" + '90.40'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'" + '  +1.80%  '
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'" + '  -0.37%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'" + '0.000008795'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'" + '  +1.84%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'" + '1.010'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'" + '  -0.27%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'" + '27.540.92'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'" + '  +1.97%  '
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'" + '  +1.27%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'" + '5.121'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'" + '  +1.12%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'" + '10.70'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'" + '  +0.43%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'" + '1.940'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'" + '  +0.53%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'" + '154.82'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'" + '  +1.24%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'" + '18.46'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'" + '  +1.33%  '
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'" + '  +0.98%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'" + '115.30'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'" + '  +0.88%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'" + '4.934'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'" + '  -0.53%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'" + '0.08884'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'" + '  -0.05%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'" + '3.321'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'" + '  -0.20%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'" + '1.214'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'" + '  +3.66%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'" + '0.7574'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'" + '  +2.13%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'" + '4.597'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'" + '  +2.03%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'" + '2.728'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'" + '  -0.39%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'" + '0.02051'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'" + '  +4.51%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'" + '1.123'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'" + '  +0.70%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'" + '0.5586'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'" + '  +7.15%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'" + '0.05276'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'" + '  -0.03%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'" + '2.990'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'" + '  +0.36%  '
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'" + '  +1.05%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'" + '8.643'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'" + '  +5.01%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'" + '0.1525'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'" + '  +0.84%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'" + '0.4914'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'" + '  +3.37%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'" + '10.67'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'" + '  +0.04%  '
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'" + '  -0.35%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'" + '1.656'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'" + '  +2.93%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'" + '102.93'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'" + '  +1.23%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'" + '67.43'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'" + '  +2.90%  '
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'" + '  +0.35%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'" + '0.9151'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'" + '  +2.94%  '
$ws.Range("E51").ClearFormats()
